$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('C3').Value = 'Open/Closed'
$ws.Range('C4').Value = 'Yes/No'
$ws.Range('C6').Value = 'CR200-034'
$ws.Range('C7').Value = 'Holbrook Justice Court'
$ws.Range('C8').Value = '03-2345'
$ws.Range('C9').Value = '06-06-2004'
$ws.Range('C10').Value = '456'
$ws.Range('B12').Value = 'Crime (Offense) Category'
$ws.Range('B14').Value = 'Case Weight'
$ws.Range('C14').Value = 'Definition? Values?'
$ws.Range('C15').Value = 'Yes/No'
$ws.Range('B20').Value = 'Case Disposition'
$ws.Range('A28').Value = 'Subject/Defendant(s)'
$ws.Range('A86').Value = 'Incidents'
$ws.Range('B86').Value = 'DR Number'
$ws.Range('B87').Value = 'Incident Date'
$ws.Range('B88').Value = 'Report Date'
$ws.Range('B89').Value = 'Location'
$ws.Range('B90').Value = 'Weapon Indicator'
$ws.Range('B91').Value = 'Arresting Agency Name'
$ws.Range('B92').Value = 'Arresting Officer'
$ws.Range('B93').Value = 'Incident Comment'
$ws.Range('A95').Value = 'Events'
$ws.Range('A96').Value = 'Offers'
$ws.Range('A97').Value = 'Research'
$ws.Range('A98').Value = 'Facts'
$ws.Range('A99').Value = 'File Location'
$ws.Range('A100').Value = 'Time'
$ws.Range('A101').Value = 'Expenses'
$ws.Range('A102').Value = 'Evidence'
$ws.Range('A103').Value = 'Related Cases'
$ws.Range('A104').Value = 'Forfeiture Cases'
$ws.Range('A105').Value = 'Case Notes'
$ws.Range('B105').Value = 'Date Created'
$ws.Range('B106').Value = 'Title'
$ws.Range('B107').Value = 'Author'
$ws.Range('B108').Value = 'Created By'
$ws.Range('A109').Value = 'File Attachment'
$ws.Range('A110').Value = 'Quick Templates'

$ws.Range("C38").Select()
